# Edit script: rebuild the "Estado de Cuenta" worker/period table with the
# updated data set (5 workers, 13 periods, new totals) and move the
# signature block down to rows 49-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Preserve the special "last row" border style (currently on row 20)
#    by copying it down to the new last data row (44) before row 20 is
#    overwritten with the regular data-row style.
# ---------------------------------------------------------------------
$lastRowStyleSrc = $ws.Range($ws.Cells.Item(20, 2), $ws.Cells.Item(20, 10))
$lastRowStyleSrc.Copy()
$newLastRow = $ws.Range($ws.Cells.Item(44, 2), $ws.Cells.Item(44, 10))
$newLastRow.PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Preserve the signature-block style (rows 25-26) by copying it down
#    to the new location (rows 49-50) before rows 25-26 are cleared.
# ---------------------------------------------------------------------
$sigLeftSrc = $ws.Range($ws.Cells.Item(25, 2), $ws.Cells.Item(26, 3))
$sigLeftSrc.Copy()
$sigLeftDst = $ws.Range($ws.Cells.Item(49, 2), $ws.Cells.Item(50, 3))
$sigLeftDst.PasteSpecial(-4122)

$sigRightSrc = $ws.Range($ws.Cells.Item(25, 8), $ws.Cells.Item(26, 10))
$sigRightSrc.Copy()
$sigRightDst = $ws.Range($ws.Cells.Item(49, 8), $ws.Cells.Item(50, 10))
$sigRightDst.PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Remove old merges (rows 25-26) and old leftover cell content, then
#    create the new merges (rows 49-50).
# ---------------------------------------------------------------------
$ws.Range("B25:C25").UnMerge()
$ws.Range("B26:C26").UnMerge()
$ws.Range("H25:J25").UnMerge()
$ws.Range("H26:J26").UnMerge()
$ws.Range("B25:J26").Clear()

$ws.Range("B49:C49").Merge()
$ws.Range("B50:C50").Merge()
$ws.Range("H49:J49").Merge()
$ws.Range("H50:J50").Merge()

# ---------------------------------------------------------------------
# 4) Paste the regular data-row style (row 16) across all data rows
#    (16-44), this also normalizes row 20 (previously special) back to
#    the regular style, and pre-applies the Text number format needed
#    for the document-number / period columns so the values we set
#    below are stored as text (matching the source data).
# ---------------------------------------------------------------------
$dataRowStyleSrc = $ws.Range($ws.Cells.Item(16, 2), $ws.Cells.Item(16, 10))
$dataRowStyleSrc.Copy()
$dataRowsDst = $ws.Range($ws.Cells.Item(16, 2), $ws.Cells.Item(43, 10))
$dataRowsDst.PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5) Fill in the worker / overdue-period table (rows 16-44).
# ---------------------------------------------------------------------
$workers = @(
    @{ doc = "45541384"; name = "LEYDIS RACERO BALASNOA"; salary = 52000; base = 1300000; periods = @("2410","2409","2408","2407") }
    @{ doc = "1002275241"; name = "EINAR ANDRES HERRERA AGAMEZ"; salary = 52000; base = 1300000; periods = @("2505","2504","2503","2502","2501","2412","2411","2410","2409","2408") }
    @{ doc = "1047404309"; name = "YESENIA DEL CARMEN MARTINEZ MARTINEZ"; salary = 56940; base = 1423500; periods = @("2506","2505") }
    @{ doc = "1091353816"; name = "MARLON EDUARDO GIL BUITRAGO"; salary = 52000; base = 1300000; periods = @("2505","2504","2503","2502","2501","2412","2411","2410","2409","2408","2407","2406") }
    @{ doc = "1003098431"; name = "DAYANA MORELO PALENCIA"; salary = 26572; base = 1423500; periods = @("2506") }
)

$row = 16
foreach ($w in $workers) {
    foreach ($p in $w.periods) {
        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $w.doc
        $ws.Cells.Item($row, 4).Value = $w.name
        $ws.Cells.Item($row, 5).Value = $p
        $ws.Cells.Item($row, 6).Value = $w.salary
        $ws.Cells.Item($row, 7).Value = $w.base
        $row = $row + 1
    }
}

# ---------------------------------------------------------------------
# 6) Update the summary header figures.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 1492452
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 13

# ---------------------------------------------------------------------
# 7) Re-fit column D (worker name) now that a longer name is present.
# ---------------------------------------------------------------------
$ws.Columns("D").AutoFit()

Write-Host "Edit complete"
